$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old row 25 ("LOB1052..." requisito line) is dropped entirely, shrinking the
# sheet from A1:C25 down to A1:C24.
$ws.Rows.Item(25).Delete()

# Cells whose content disappears in the new layout (fully cleared, not just blanked)
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()

# Re-populate row 10 and rows 13-24 with the rearranged content
$ws.Range("B10").Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Range("C10").Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A15").Value = "Programa:"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Range("C18").Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("C19").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("C20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOB1019 -  Física II  (Requisito)
"
$ws.Range("C23").Value = "LOB1019 -  Física II  (Requisito)
"
$ws.Range("B24").Value = "LOB1052 -  Cálculo III  (Requisito)
"
$ws.Range("C24").Value = "LOB1052 -  Cálculo III  (Requisito)
"

# B15/C15 need the literal text "01/01/2012" (not an actual date). Copying the
# value from B8/C8 (which already hold this same text) keeps it text instead of
# Excel auto-converting a typed "01/01/2012" into a date serial; the format is
# then copied back from B9/C9 so the cells keep their normal (non-date) style.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("B9").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("C9").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# B18 and B23 are brand-new cells (column B had nothing in these rows before), so
# they picked up the neighbouring column A bold style; copy column Bs normal
# (wrap-text) format back onto them from another row that already has it.
$ws.Range("B9").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights for the new layout
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30

# Rows 17 and 22 revert to the default (non-custom) row height
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(22).AutoFit()
